{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// [rowIndex, columnIndex, newText] updates to the division-problem cells,\n// applied in document (top-to-bottom, left-to-right) order.\nconst updates = [\n  [0, 0, \"33\u00f73=\"],\n  [0, 1, \"12\u00f76=\"],\n  [0, 2, \"53\u00f76=\"],\n  [0, 3, \"66\u00f76=\"],\n  [0, 4, \"36\u00f77=\"],\n  [4, 0, \"31\u00f77=\"],\n  [4, 1, \"67\u00f77=\"],\n  [4, 2, \"47\u00f74=\"],\n  [4, 3, \"49\u00f73=\"],\n  [4, 4, \"74\u00f73=\"],\n  [8, 0, \"25\u00f72=\"],\n  [8, 1, \"35\u00f72=\"],\n  [8, 2, \"68\u00f77=\"],\n  [8, 3, \"59\u00f73=\"],\n  [8, 4, \"59\u00f74=\"],\n  [12, 0, \"53\u00f78=\"],\n  [12, 1, \"81\u00f78=\"],\n  [12, 2, \"43\u00f77=\"],\n  [12, 3, \"90\u00f77=\"],\n  [12, 4, \"45\u00f76=\"],\n  [16, 0, \"61\u00f72=\"],\n  [16, 1, \"77\u00f73=\"],\n  [16, 2, \"57\u00f78=\"],\n  [16, 3, \"79\u00f79=\"],\n  [16, 4, \"83\u00f74=\"],\n];\n\nfor (const [rowIndex, columnIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, columnIndex);\n  const range = cell.body.getRange(\"Whole\");\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row/Col (1-based) -> new cell text, applied in document (top-to-bottom,\n# left-to-right) order, matching the division-problem cells only.\n$updates = @(\n    @{Row = 1; Col = 1; Text = \"33\u00f73=\"}\n    @{Row = 1; Col = 2; Text = \"12\u00f76=\"}\n    @{Row = 1; Col = 3; Text = \"53\u00f76=\"}\n    @{Row = 1; Col = 4; Text = \"66\u00f76=\"}\n    @{Row = 1; Col = 5; Text = \"36\u00f77=\"}\n    @{Row = 5; Col = 1; Text = \"31\u00f77=\"}\n    @{Row = 5; Col = 2; Text = \"67\u00f77=\"}\n    @{Row = 5; Col = 3; Text = \"47\u00f74=\"}\n    @{Row = 5; Col = 4; Text = \"49\u00f73=\"}\n    @{Row = 5; Col = 5; Text = \"74\u00f73=\"}\n    @{Row = 9; Col = 1; Text = \"25\u00f72=\"}\n    @{Row = 9; Col = 2; Text = \"35\u00f72=\"}\n    @{Row = 9; Col = 3; Text = \"68\u00f77=\"}\n    @{Row = 9; Col = 4; Text = \"59\u00f73=\"}\n    @{Row = 9; Col = 5; Text = \"59\u00f74=\"}\n    @{Row = 13; Col = 1; Text = \"53\u00f78=\"}\n    @{Row = 13; Col = 2; Text = \"81\u00f78=\"}\n    @{Row = 13; Col = 3; Text = \"43\u00f77=\"}\n    @{Row = 13; Col = 4; Text = \"90\u00f77=\"}\n    @{Row = 13; Col = 5; Text = \"45\u00f76=\"}\n    @{Row = 17; Col = 1; Text = \"61\u00f72=\"}\n    @{Row = 17; Col = 2; Text = \"77\u00f73=\"}\n    @{Row = 17; Col = 3; Text = \"57\u00f78=\"}\n    @{Row = 17; Col = 4; Text = \"79\u00f79=\"}\n    @{Row = 17; Col = 5; Text = \"83\u00f74=\"}\n)\n\nforeach ($u in $updates) {\n    $table.Cell($u.Row, $u.Col).Range.Text = $u.Text\n}\n"}
